# -----------------------------------------------------------------------
# Example_Table.xlsx update:
#   - TitleBlockData sheet: table grows from A1:E17 to A1:E21, column
#     "Multiplier" renamed to "Factor", property rows re-sorted
#     alphabetically with several new / changed values, column widths
#     adjusted.
#   - Settings sheet: DrwNrFieldName now points at "DN" and MapMass is
#     cleared.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TitleBlockData")

# --- materialize the new rows (19-22) with the same formatting used by
#     the rest of the table (style carries left/vcenter/indent alignment) --
$ws.Range("A17:E17").Copy()
$ws.Range("A18:E22").PasteSpecial(-4122)  # xlPasteFormats

# --- resize the table / autofilter to cover the extra rows -------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E21"))

# --- header row: "Multiplier" -> "Factor" -------------------------------
$ws.Range("D1").Value = "Factor"

# --- column widths -------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 18 - 0.8333333333333333
$ws.Columns.Item(2).ColumnWidth = 84 - 0.8333333333333333
$ws.Columns.Item(4).ColumnWidth = 11 - 0.8333333333333333

# --- property rows (alphabetically sorted Property Name column) --------
$rows = @(
    @(2, "ANGLE UNIT", "°", $null),
    @(3, "APPROVER_NAME", "APPROVER NAME", $null),
    @(4, "AUTHOR_NAME", $null, $null),
    @(5, "DATE", $null, $null),
    @(6, "DN", "Test object #2", $null),
    @(7, "DOCUMENT_TYPE", "Mechanical assembly drawing", $null),
    @(8, "LENGTH UNIT", $null, $null),
    @(9, "OWNER_NAME", "FreeCAD", $null),
    @(10, "PM", $null, $null),
    @(11, "PN", $null, $null),
    @(12, "PROJECT NAME", "Titleblock Workbench", $null),
    @(13, "REVISION", "REV A", $null),
    @(14, "RIGHTS", "(R) DO NOT DUPLICATE THIS DRAWING TO THIRD PARTIES WITHOUT OWNER'S PERMISSION !", $null),
    @(15, "ROUGHNESS", "3.2", $null),
    @(16, "SCALE", "M x:x", $null),
    @(17, "SHEET", "1", "X"),
    @(18, "SIZE", "A3", $null),
    @(19, "TITLELINE-1", "Title", $null),
    @(20, "TITLELINE-2", "Module name", $null),
    @(21, "TOLERANCE", "+/- ?", $null)
)

# cell refs whose value looks like a number but must stay text (matches
# the source workbook, which stores every property value as a string)
$forceTextCells = @("B15", "B17")

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]

    $bRef = "B$r"
    if (($forceTextCells -contains $bRef) -and ($row[2] -ne $null)) {
        # A leading apostrophe forces a text entry instead of a number;
        # then re-apply the shared formatting (copied from a plain
        # style-1 cell) so the cell keeps the same style index as the
        # rest of the table instead of Excel's auto "quote prefix" style
        $ws.Range($bRef).Value = "'" + $row[2]
        $ws.Range("A16").Copy()
        $ws.Range($bRef).PasteSpecial(-4122)
    } else {
        $ws.Range($bRef).Value = $row[2]
    }

    $ws.Range("C$r").Value = $row[3]
}

# --- Settings sheet bug fixes -------------------------------------------
$ws2 = $wb.Worksheets.Item("Settings")
$ws2.Range("B11").Value = "DN"
$ws2.Range("B14").Value = $null
